$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2561")
$ws.Range("F26").Formula = "12235.8"
